$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are always stored as text so that values such as
# "608.69" or "3.147.73" are preserved exactly instead of being coerced to numbers.

# Update Price (D) and Volume(1h) (E) columns for rows with simple value changes
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.184.70"
$ws.Range("E2").Value = "  +5.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.147.73"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.88"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "608.69"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.11"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.381"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.147.67"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.785"
$ws.Range("E11").Value = "  -5.00%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.627.92"
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.730.45"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.141.32"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "523.64"
$ws.Range("E19").Value = "  +18.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.39"
$ws.Range("E20").Value = "  -7.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.48"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  -5.40%  "
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.70"
$ws.Range("E24").Value = "  -3.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "89.96"
$ws.Range("E25").Value = "  +5.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.43"
$ws.Range("E26").Value = "  -4.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.50"
$ws.Range("E27").Value = "  -9.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.310.70"
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  -3.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.174"
$ws.Range("E31").Value = "  -4.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "8.89"
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.66"
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.151"
$ws.Range("E36").Value = "  -5.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.19"
$ws.Range("E37").Value = "  -9.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.37"
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.433"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "464.73"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.21"
$ws.Range("E42").Value = "  -6.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.51"
$ws.Range("E43").Value = "  -12.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.09"
$ws.Range("E45").Value = "  -5.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "162.67"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.89"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.693"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.46"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.18"
$ws.Range("E50").Value = "  +0.81%  "

# Rows where the coin position/ranking swapped or changed entirely (full row content rewrite)
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.82"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.776"
$ws.Range("E51").Value = "  +6.73%  "
